# Tweak wording, add test cases and update storage diagram
# (this script only touches the parts relevant to this presentation file)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" date placeholder text from
#    11/10/2018 -> 11/12/2018 everywhere it appears: every slide layout,
#    the slide master, and the notes master.
# ---------------------------------------------------------------------
$ppPlaceholderDate = 16
$oldDate = "11/10/2018"
$newDate = "11/12/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $ph = $shp.PlaceholderFormat
            if ($ph.Type -eq $ppPlaceholderDate) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Every custom layout hanging off the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder($layouts.Item($li).Shapes)
}

# The slide master itself.
Update-DatePlaceholder($p.SlideMaster.Shapes)

# The notes master.
Update-DatePlaceholder($p.NotesMaster.Shapes)

# ---------------------------------------------------------------------
# 2) Rename the "AddressBook" class box on slide 1 to "ModulePlanner".
#    Use a targeted Characters() sub-range so only the second run's text
#    is touched (keeps the existing line break / run formatting intact).
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq 66 -and $shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $fullText = $tr.Text
        $pos = $fullText.IndexOf("AddressBook")
        if ($pos -ge 0) {
            $sub = $tr.Characters($pos + 1, "AddressBook".Length)
            $sub.Text = "ModulePlanner"
        }
    }
}
